$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 11.1339517624806
$ws.Range("C2").Value = 9.19092563734667
$ws.Range("D2").Value = 3.62315315705178
$ws.Range("F2").Value = 17.10648339475095
$ws.Range("G2").Value = 16.53636565325821
$ws.Range("H2").Value = 11.11538305794869
$ws.Range("I2").Value = 15.39394819057766
$ws.Range("M2").Value = 19.96649101233919
$ws.Range("O2").Value = 15.27676452931069
$ws.Range("B3").Value = 10.51735665172442
$ws.Range("C3").Value = 8.887724939487258
$ws.Range("D3").Value = 3.519532704283068
$ws.Range("F3").Value = 17.11250326720325
$ws.Range("G3").Value = 16.51477894963819
$ws.Range("H3").Value = 11.16225475457281
$ws.Range("I3").Value = 15.51713480260586
$ws.Range("M3").Value = 19.35727919110772
$ws.Range("O3").Value = 15.34255989920806
$ws.Range("B4").Value = 10.11788713663211
$ws.Range("C4").Value = 8.695476398131865
$ws.Range("D4").Value = 3.453700092730103
$ws.Range("F4").Value = 17.1228064795536
$ws.Range("G4").Value = 16.51127774048831
$ws.Range("H4").Value = 11.19330796370524
$ws.Range("I4").Value = 15.59701654360547
$ws.Range("M4").Value = 18.98144038684232
$ws.Range("O4").Value = 15.38773332158033
$ws.Range("B5").Value = 9.949936088562673
$ws.Range("C5").Value = 8.615698730841677
$ws.Range("D5").Value = 3.426342123454693
$ws.Range("F5").Value = 17.12866217513389
$ws.Range("G5").Value = 16.51229570105033
$ws.Range("H5").Value = 11.20653343761598
$ws.Range("I5").Value = 15.63063649100059
$ws.Range("M5").Value = 18.82809404978684
$ws.Range("O5").Value = 15.40733671287853
$ws.Range("B6").Value = 9.921738948799137
$ws.Range("C6").Value = 8.60236799522777
$ws.Range("D6").Value = 3.421768072974222
$ws.Range("F6").Value = 17.12973443472573
$ws.Range("G6").Value = 16.51261207757898
$ws.Range("H6").Value = 11.20876398572782
$ws.Range("I6").Value = 15.63628353702683
$ws.Range("M6").Value = 18.80262741194279
$ws.Range("O6").Value = 15.41066383448011
$ws.Range("B7").Value = 10.11564287122398
$ws.Range("C7").Value = 8.694406166484567
$ws.Range("D7").Value = 3.453333247721986
$ws.Range("F7").Value = 17.12287874904364
$ws.Range("G7").Value = 16.51128158392738
$ws.Range("H7").Value = 11.19348401593062
$ws.Range("I7").Value = 15.59746563155766
$ws.Range("M7").Value = 18.97937271818323
$ws.Range("O7").Value = 15.38799286967773
$ws.Range("B8").Value = 10.9257286892898
$ws.Range("C8").Value = 9.087698265458195
$ws.Range("D8").Value = 3.58789762012811
$ws.Range("F8").Value = 17.10718562348962
$ws.Range("G8").Value = 16.52689495306661
$ws.Range("H8").Value = 11.13107205627235
$ws.Range("I8").Value = 15.43554202046106
$ws.Range("M8").Value = 19.75695518406187
$ws.Range("O8").Value = 15.29845632054804
$ws.Range("B9").Value = 12.34576091478683
$ws.Range("C9").Value = 9.806990491616837
$ws.Range("D9").Value = 3.83326349442646
$ws.Range("F9").Value = 17.12895677297172
$ws.Range("G9").Value = 16.63503012637318
$ws.Range("H9").Value = 11.02675633285768
$ws.Range("I9").Value = 15.15168428666479
$ws.Range("M9").Value = 21.25704797118363
$ws.Range("O9").Value = 15.16101361621645
$ws.Range("B10").Value = 13.2834387030545
$ws.Range("C10").Value = 10.29944812860641
$ws.Range("D10").Value = 4.001073446550656
$ws.Range("F10").Value = 17.17704619437959
$ws.Range("G10").Value = 16.76159933680144
$ws.Range("H10").Value = 10.96117814056006
$ws.Range("I10").Value = 14.96365628838721
$ws.Range("M10").Value = 22.3310365653944
$ws.Range("O10").Value = 15.08361813407585
$ws.Range("B11").Value = 13.68669570836184
$ws.Range("C11").Value = 10.514907170276
$ws.Range("D11").Value = 4.074499827138403
$ws.Range("F11").Value = 17.20587066282004
$ws.Range("G11").Value = 16.82928727925874
$ws.Range("H11").Value = 10.93375670384767
$ws.Range("I11").Value = 14.88257461944381
$ws.Range("M11").Value = 22.81111071181895
$ws.Range("O11").Value = 15.05359713305101
$ws.Range("B12").Value = 13.83602693235236
$ws.Range("C12").Value = 10.59521147537444
$ws.Range("D12").Value = 4.10187077952157
$ws.Range("F12").Value = 17.21778035260748
$ws.Range("G12").Value = 16.85635581667532
$ws.Range("H12").Value = 10.92372038477551
$ws.Range("I12").Value = 14.8525119645525
$ws.Range("M12").Value = 22.99150595583847
$ws.Range("O12").Value = 15.0429799612694
$ws.Range("B13").Value = 13.80401611795768
$ws.Range("C13").Value = 10.57797446417964
$ws.Range("D13").Value = 4.095995485961135
$ws.Range("F13").Value = 17.2151712501355
$ws.Range("G13").Value = 16.85046255209955
$ws.Range("H13").Value = 10.92586641128858
$ws.Range("I13").Value = 14.85895796206662
$ws.Range("M13").Value = 22.95271944709279
$ws.Range("O13").Value = 15.04523307081666
$ws.Range("B14").Value = 13.69904899711841
$ws.Range("C14").Value = 10.52153990058702
$ws.Range("D14").Value = 4.076760432607545
$ws.Range("F14").Value = 17.20683058193724
$ws.Range("G14").Value = 16.8314855474578
$ws.Range("H14").Value = 10.93292403817213
$ws.Range("I14").Value = 14.88008849052363
$ws.Range("M14").Value = 22.82598091785019
$ws.Range("O14").Value = 15.0527085687339
$ws.Range("B15").Value = 13.63431370323836
$ws.Range("C15").Value = 10.48680321728668
$ws.Range("D15").Value = 4.064921455901983
$ws.Range("F15").Value = 17.20185102593189
$ws.Range("G15").Value = 16.82004808148557
$ws.Range("H15").Value = 10.93729233294593
$ws.Range("I15").Value = 14.89311508172285
$ws.Range("M15").Value = 22.74816282512138
$ws.Range("O15").Value = 15.0573854992693
$ws.Range("B16").Value = 13.25661395792485
$ws.Range("C16").Value = 10.28519029853633
$ws.Range("D16").Value = 3.996214924247871
$ws.Range("F16").Value = 17.17530190535094
$ws.Range("G16").Value = 16.7573778430197
$ws.Range("H16").Value = 10.9630187823774
$ws.Range("I16").Value = 14.96904484543353
$ws.Range("M16").Value = 22.29947617078089
$ws.Range("O16").Value = 15.08568490392417
$ws.Range("B17").Value = 13.01892158867771
$ws.Range("C17").Value = 10.15927362576457
$ws.Range("D17").Value = 3.953308567450247
$ws.Range("F17").Value = 17.16079173317798
$ws.Range("G17").Value = 16.72151036063241
$ws.Range("H17").Value = 10.97941917862786
$ws.Range("I17").Value = 15.01676642995571
$ws.Range("M17").Value = 22.02191508341014
$ws.Range("O17").Value = 15.10437793826157
$ws.Range("B18").Value = 12.88001528174837
$ws.Range("C18").Value = 10.08604653767604
$ws.Range("D18").Value = 3.928356758000629
$ws.Range("F18").Value = 17.15310026605612
$ws.Range("G18").Value = 16.70183334512755
$ws.Range("H18").Value = 10.98907907395778
$ws.Range("I18").Value = 15.04463372158566
$ws.Range("M18").Value = 21.86147976371106
$ws.Range("O18").Value = 15.1156175625669
$ws.Range("B19").Value = 12.83260834827503
$ws.Range("C19").Value = 10.06111692114148
$ws.Range("D19").Value = 3.919862077254167
$ws.Range("F19").Value = 17.15060857578843
$ws.Range("G19").Value = 16.69533518485142
$ws.Range("H19").Value = 10.9923886757888
$ws.Range("I19").Value = 15.05414106460182
$ws.Range("M19").Value = 21.80702923172069
$ws.Range("O19").Value = 15.11950673511267
$ws.Range("B20").Value = 13.04445139527043
$ws.Range("C20").Value = 10.17276120640521
$ws.Range("D20").Value = 3.957904415223175
$ws.Range("F20").Value = 17.16226866632512
$ws.Range("G20").Value = 16.72522999444498
$ws.Range("H20").Value = 10.97764984631551
$ws.Range("I20").Value = 15.01164300503403
$ws.Range("M20").Value = 22.05154499251239
$ws.Range("O20").Value = 15.10233749932851
$ws.Range("B21").Value = 13.72997209252052
$ws.Range("C21").Value = 10.5381513824209
$ws.Range("D21").Value = 4.082422127752193
$ws.Range("F21").Value = 17.2092534938911
$ws.Range("G21").Value = 16.83702072477289
$ws.Range("H21").Value = 10.93084160100243
$ws.Range("I21").Value = 14.87386453232088
$ws.Range("M21").Value = 22.86324637638478
$ws.Range("O21").Value = 15.05049240737029
$ws.Range("B22").Value = 14.15833004832103
$ws.Range("C22").Value = 10.76944239748991
$ws.Range("D22").Value = 4.161265838082515
$ws.Range("F22").Value = 17.245754076089
$ws.Range("G22").Value = 16.91844491694499
$ws.Range("H22").Value = 10.9022761441352
$ws.Range("I22").Value = 14.78755630993841
$ws.Range("M22").Value = 23.38551920272064
$ws.Range("O22").Value = 15.02098907330818
$ws.Range("B23").Value = 13.93151314807053
$ws.Range("C23").Value = 10.64670161777727
$ws.Range("D23").Value = 4.119422152963404
$ws.Range("F23").Value = 17.22574492482722
$ws.Range("G23").Value = 16.87422895122889
$ws.Range("H23").Value = 10.91733632555325
$ws.Range("I23").Value = 14.83327829180767
$ws.Range("M23").Value = 23.1075786994847
$ws.Range("O23").Value = 15.03633306452442
$ws.Range("B24").Value = 13.03291638794936
$ws.Range("C24").Value = 10.16666607071357
$ws.Range("D24").Value = 3.95582751701864
$ws.Range("F24").Value = 17.16159891796611
$ws.Range("G24").Value = 16.72354540738666
$ws.Range("H24").Value = 10.97844904173639
$ws.Range("I24").Value = 15.01395796142437
$ws.Range("M24").Value = 22.03815198109268
$ws.Range("O24").Value = 15.10325844736921
$ws.Range("B25").Value = 11.97999144143479
$ws.Range("C25").Value = 9.618444587814226
$ws.Range("D25").Value = 3.768996502640696
$ws.Range("F25").Value = 17.11742893297965
$ws.Range("G25").Value = 16.59747399047141
$ws.Range("H25").Value = 11.05303701956284
$ws.Range("I25").Value = 15.22487138325943
$ws.Range("M25").Value = 20.85529671751131
$ws.Range("O25").Value = 15.19407731386453
